# Refresh the crypto price/volume table (cryptos.xlsx) with the
# latest scrape: updates Price (D) and Volume(1h) (E) for every
# coin row, and swaps the EOS/Quant rows (49/50) back in rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.355.31"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.865.87"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").Value = "'1.020"
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'316.57"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'1.018"
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").Value = "'0.5098"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").Value = "'0.3955"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "'0.08348"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "'1.107"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'41.90"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'6.226"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").Value = "1.830.59"
$ws.Range("E14").Value = "  -4.26%  "
$ws.Range("D15").Value = "'1.019"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'7.194"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "'0.00001105"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'90.42"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "'0.06735"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "'1.018"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'5.954"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").Value = "28.363.00"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'2.288"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'161.79"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "2.038.35"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D29").Value = "'2.355"
$ws.Range("E29").Value = "  -5.18%  "
$ws.Range("D30").Value = "'127.05"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("D31").Value = "'0.1045"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'1.030"
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("D33").Value = "'5.772"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'3.641"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.02418"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").Value = "'0.06465"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  -1.72%  "
$ws.Range("D38").Value = "'8.832"
$ws.Range("E38").Value = "  -8.52%  "
$ws.Range("D39").Value = "'1.264"
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "'0.6367"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("D42").Value = "'4.988"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'0.6013"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "'13.04"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'3.707"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "'1.213"
$ws.Range("E47").Value = "  -5.76%  "
$ws.Range("D48").Value = "'1.983"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "'1.200"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'120.83"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'0.06840"
$ws.Range("E51").Value = "  -1.44%  "
